# Applies the commit "Add files via upload":
#   - Removes the Solution 1-6 write-up text from the Dataset sheet
#     (cells A29:A31, A33:A35, A37:A40, A42:A45, A47:A49, A51:A53, A55:A58 -
#     the blank spacer rows A32/A36/A41/A46/A50/A54 are left untouched),
#     which in turn shrinks sharedStrings.xml (the 24 strings become unused).
#   - Updates the Dataset sheet's active selection to D26 (and drops the
#     A39 scroll position in the process).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

$ws.Activate()

$cellsToClear = @(
    "A29","A30","A31",
    "A33","A34","A35",
    "A37","A38","A39","A40",
    "A42","A43","A44","A45",
    "A47","A48","A49",
    "A51","A52","A53",
    "A55","A56","A57","A58"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Match the author's final selection / cursor position on the sheet.
$ws.Range("D26").Select()
